# Update "想去人数" (want-to-go count) values in "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$updates1 = @{
    "F4" = 264
    "F9" = 365
    "F11" = 433
    "F13" = 2549
    "F15" = 1346
    "F16" = 4702
    "F18" = 5116
    "F19" = 1672
    "F20" = 2869
    "F21" = 3267
    "F22" = 168
    "F23" = 1557
    "F24" = 258
    "F25" = 837
    "F26" = 108
    "F27" = 293
    "F28" = 1000
    "F29" = 1871
    "F31" = 281
    "F32" = 709
    "F33" = 158
    "F35" = 415
}

$updates2 = @{
    "F10" = 264
    "F14" = 365
    "F18" = 433
    "F21" = 2549
    "F22" = 1346
    "F26" = 4702
    "F28" = 5116
    "F29" = 1672
    "F30" = 2869
    "F31" = 3267
    "F32" = 168
    "F35" = 1557
    "F37" = 258
    "F38" = 837
    "F39" = 108
    "F40" = 293
    "F41" = 1000
    "F43" = 1871
    "F45" = 281
    "F46" = 709
    "F47" = 158
    "F49" = 415
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $updates1.Keys) {
    $ws1.Range($addr).Value = $updates1[$addr]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $updates2.Keys) {
    $ws4.Range($addr).Value = $updates2[$addr]
}

Write-Host "Update complete."
